$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
  '3535 Opal Meadow Heights Aged Care Community Meadow Heights',
  '7 Chefs Fawkner',
  'A1 Bakery Brunswick',
  'Acquire BPO Southbank',
  'Al Haj Halal Meats Glenroy',
  'Al-Taqwa College Truganina',
  'Amiga Montessori Craigieburn',
  'Best&Less Fountain Gate Narre Warren',
  'Broadmeadows Medical Centre Broadmeadows',
  'Budget Car and Truck Rentals Campbellfield',
  'Campbellfield Heights Primary School Campbellfield',
  'Can Panel Campbellfield',
  'Cannie Road Construction Site Cannie',
  'Caroline Springs Police Station',
  'Cedars Medical Clinic Coburg',
  'Chemist Warehouse Fillo Drive Somerton',
  'City of Hobsons Bay Community',
  'City of Moreland Community',
  'Classy Cabinets and Kitchens Craigieburn',
  'Coles Aurora Village Epping',
  'Coles Barkly Square Brunswick August',
  'Coles Broadmeadows Central Shopping Centre',
  'Coles Campbellfield Plaza Campbellfield',
  'Coles Coburg North Village',
  'Coles Coburg North Village Aug',
  'Coles Greenvale Shopping Centre',
  'Coles Pakenham Place Shopping Centre',
  'Coles Roxburgh Village Roxburgh Park',
  'Community Kids Meadow Heights',
  'Construction Site Olea Apartment Caulfield North',
  'Costco Wholesale Epping',
  'Croydon Orthodontics',
  'Direct Freight Express Campbellield',
  'Elite Smart Community Care Campbellfield',
  'Ernst and Young 8 Exhibition Street Melbourne',
  'Fitzroy Community School Fitzroy North',
  'Glenroy West Primary School',
  'Hamilton Marino 236 Jasper Road McKinnon',
  'Health Care Providers Association South Melbourne',
  'IGA Meadow Heights Shopping Centre Meadow Heights',
  'Ilim College Glenroy Campus Hadfield',
  'Ilim College Kiewa Campus Boys Secondary Dallas',
  'Ilim Leaning Sanctuary Glenroy',
  'Industrial Galvanizers Valmont Coatings Campbellfield',
  'Islamic College of Melbourne Tarneit',
  'KFC Fawkner',
  'Kasr Sweets Coolaroo',
  'Kids House Early Learning Cheltenham',
  'Learning Nest Early Learning Centre Meadow Heights',
  'Level Crossing Removal Project Lilydale Construction Site John Street',
  'Lineage Logistics Laverton North',
  'Linfox Somerton National Distribution Centre Somerton',
  'Malvern Health and Fitness Clinic Malvern',
  'McDonalds Thomastown II',
  'Melbourne Metropolitan Remand Centre Ravenhall',
  'Melbourne Truck Repairs Campbellfield',
  'Melbourne West Police Station Docklands',
  'Melbourne Youth Justice Centre Parkville',
  'Melton Police Station Melton',
  'Mercy Hospital for Women Heidelberg',
  'Mernda YMCA Early Learning Centre Mernda',
  'Montessori Beginnings Greenvale',
  'MyCentre Childcare Broadmeadows',
  'Newbury Child and Community Centre Craigieburn',
  'Newport Football Club Altona North',
  'Newport Gardens Early Years Centre Newport',
  'Nido Early School Moonee Ponds',
  'Nino Early Learning Adventures Lalor',
  'North Geelong House Party',
  'Northern Health Northern Hospital Epping Emergency Department Tier 1B',
  'OnQ Plumbing and Excavations Craigieburn',
  'Oporto Coolaroo',
  'Paisley Park Early Learning Centre Bundoora',
  'Panorama Construction Site Whitehorse Rd Box Hill',
  'People First Healthcare Home Residence Disability Support Taylors Lakes',
  'Ramsay Health Care Warringal Private Hospital Heidelberg',
  'Salta Drive Construction Site Rangedale Drainage Altona North',
  'Serco Mill Park',
  'St Vincents Hospital Emergency Department Melbourne',
  'Tek Foods Somerton',
  'The Homestead Child and Family Centre Roxburgh Park',
  'The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B',
  'Tip Top Warehouse Dandenong',
  'Tunstall Fresh Tunstall Square Shopping Centre Doncaster East',
  'Unilodge College Square Student Accommodation 570 Lygon Street Carlton',
  'Werribee Mercy Hospistal Emergency Department',
  'Western Health Sunshine Hospital Emergency Department',
  'Who is Bunker Spreckels Cafe Elwood',
  'Woodlands Long Day Care and Kindergarten Roxburgh Park',
  'Woolworths Greenvale Lakes Roxburgh Park',
  'Yarra Childcare Centre Truganina'
)

$values = @(
  22,
  6,
  5,
  12,
  53,
  12,
  25,
  7,
  5,
  7,
  6,
  5,
  7,
  8,
  44,
  5,
  10,
  6,
  13,
  5,
  5,
  10,
  8,
  15,
  7,
  6,
  7,
  11,
  19,
  8,
  21,
  5,
  6,
  5,
  5,
  55,
  7,
  6,
  16,
  6,
  22,
  5,
  15,
  15,
  8,
  11,
  9,
  13,
  9,
  8,
  7,
  6,
  5,
  6,
  9,
  6,
  6,
  5,
  5,
  6,
  5,
  6,
  15,
  10,
  5,
  5,
  18,
  5,
  7,
  52,
  15,
  17,
  10,
  28,
  5,
  9,
  6,
  5,
  6,
  11,
  18,
  7,
  10,
  5,
  13,
  7,
  7,
  5,
  5,
  10,
  15
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}